# Phase 3 RAD Non-UI Test Cases and Data
# Appends the newly-added TaxType test rows (rows 30-50) to Sheet1 and
# updates the active selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Y", "Existing Liability w/Notice Number", "Admissions and Amusement Tax"),
    @("Y", "Existing Liability w/Notice Number", "Estate Tax"),
    @("Y", "Existing Liability w/Notice Number", "Motor Fuel Tax"),
    @("Y", "Existing Liability w/Notice Number", "Slots License Fee"),
    @("Y", "Existing Liability w/Notice Number", "Tobacco Tax"),
    @("Y", "Existing Liability w/Notice Number", "Transportation Network Services"),
    @("Y", "Existing Liability w/Notice Number", "Unclaimed Property"),
    @("Y", "Existing Liability w/Notice Number", "IFTA Tax"),
    @("Y", "New Tax Return Amount Due", "Admissions and Amusement Tax"),
    @("Y", "New Tax Return Amount Due", "Alcohol Tax"),
    @("Y", "New Tax Return Amount Due", "Bay Restoration Fee"),
    @("Y", "New Tax Return Amount Due", "Corporate Income Tax"),
    @("Y", "New Tax Return Amount Due", "Estate Tax"),
    @("Y", "New Tax Return Amount Due", "Motor Fuel Tax"),
    @("Y", "New Tax Return Amount Due", "Sales and Use Tax"),
    @("Y", "New Tax Return Amount Due", "Slots License Fee"),
    @("Y", "New Tax Return Amount Due", "Tire Recycling Fee"),
    @("Y", "New Tax Return Amount Due", "Tobacco Tax"),
    @("Y", "New Tax Return Amount Due", "Transportation Network Services"),
    @("Y", "New Tax Return Amount Due", "Unclaimed Property"),
    @("Y", "New Tax Return Amount Due", "Withholding Tax")
)

$startRow = 30
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 3).Value = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $r = $r + 1
}

$lastRow = $r - 1

# Match the workbook's new selection / scroll position.
$ws.Range("C25:C" + $lastRow).Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 29
$win.ScrollColumn = 1
